$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.869.42'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.625.22'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.49'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.854.89'
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').Value = '1.618.01'
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '27.859.21'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.88'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').Value = '0.0₃0722'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.22%  '
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('D34').Value = '1.397.17'
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E36').Value = '  +9.23%  '
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.861'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.82'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').Value = '1.764.86'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0103'
$ws.Range('E50').Value = '  -2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0504'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.43%  '
